$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "34×30=1020" "66×54=3564"
Replace-Text "99×30=2970" "70×52=3640"
Replace-Text "77×98=7546" "48×68=3264"
Replace-Text "46×31=1426" "60×57=3420"
Replace-Text "23×87=2001" "88×94=8272"
Replace-Text "36×98=3528" "67×71=4757"
Replace-Text "11×20=220" "75×40=3000"
Replace-Text "67×51=3417" "67×89=5963"
Replace-Text "76×37=2812" "31×64=1984"
Replace-Text "36×59=2124" "53×31=1643"
Replace-Text "51×78=3978" "80×62=4960"
Replace-Text "48×70=3360" "37×32=1184"
Replace-Text "83×49=4067" "50×18=900"
Replace-Text "96×16=1536" "13×30=390"
Replace-Text "23×11=253" "30×17=510"
Replace-Text "56×30=1680" "54×19=1026"
Replace-Text "22×86=1892" "74×24=1776"
Replace-Text "47×30=1410" "44×69=3036"
Replace-Text "58×84=4872" "56×39=2184"
Replace-Text "58×99=5742" "24×24=576"
Replace-Text "87×61=5307" "19×24=456"
Replace-Text "77×73=5621" "20×99=1980"
Replace-Text "61×55=3355" "20×74=1480"
Replace-Text "76×47=3572" "32×22=704"
Replace-Text "28×96=2688" "30×71=2130"
